$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data rows to append (rows 10, 11, 12)
$data = @(
    @{ Row = 10; A = 8914.68;            B = 8995.64;            C = 18.84; D = 19.010000000000002;  E = $true;  F = 0.9;                    G = 42613.765532407408;  H = $false },
    @{ Row = 11; A = 8963.7099999999991;  B = 8914.68;            C = 18.93; D = 18.824999999999999;  E = $true;  F = -0.55000000000000004;  G = 42614.672847222224;  H = $true },
    @{ Row = 12; A = 9011.2199999999993;  B = 8963.7099999999991; C = 18.72; D = 18.62;                E = $true;  F = -0.53;                  G = 42615.750138888892;  H = $true }
)

foreach ($d in $data) {
    $r = $d.Row
    $ws.Cells.Item($r, 1).Value = $d.A
    $ws.Cells.Item($r, 2).Value = $d.B
    $ws.Cells.Item($r, 3).Value = $d.C
    $ws.Cells.Item($r, 4).Value = $d.D
    $ws.Cells.Item($r, 5).Value = $d.E
    $ws.Cells.Item($r, 6).Value = $d.F
    $ws.Cells.Item($r, 7).Value = $d.G
    $ws.Cells.Item($r, 7).NumberFormat = "m/d/yy h:mm"
    $ws.Cells.Item($r, 8).Value = $d.H
}
